$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 2")

# --- Insert two new rows above the "Estimate Totals" row (old row 26) ---
# This shifts the totals row (26 -> 28) and the bug-note rows (28,29 -> 30,31)
# down, and inherits formatting from the row above (row 25), matching the
# target styles for the new rows.
$ws.Rows("26:27").Insert()

# --- Populate the two new rows (26 and 27) ---
$ws.Range("A26").Value = "Inventory"
$ws.Range("C26").Value = "Adjust Component Add/Edit to only allow valid inputs (No numbers on costs etc)"
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 0
$ws.Range("H26").Value = 0

$ws.Range("A27").Value = "Inventory"
$ws.Range("C27").Value = "Adjust Product Add/Edit to only allow valid inputs (No numbers on costs etc)"
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 0

# --- Update the chart's plotted range to follow the totals row's new location ---
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$ser = $chart.SeriesCollection(1)
$ser.Formula = "=SERIES(,,'Sprint 2'!`$D`$28:`$H`$28,1)"

# --- Widen column C to fit the new, longer text ---
# (ColumnWidth uses Excel's "characters" unit, which is offset from the raw
# OOXML <col width> by the default-font padding; 73.17087140324561 is the
# ColumnWidth value that round-trips to an OOXML width of exactly 74.)
$ws.Columns("C").ColumnWidth = 73.17087140324561

# --- Update the view: scroll position, zoom, and active selection ---
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Application.ActiveWindow.Zoom = 160
$ws.Range("I28").Select()
